# Update cached market-price / profit figures across the Goblin Profits
# workbook (scheduled market-data refresh). Each worksheet corresponds to
# a crafting class; columns H-N hold the price/profit figures that were
# recalculated from refreshed market data. A handful of rows gained or
# lost an H-N cell entirely (no HQ/NQ price recorded that refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 931.36365
$ws.Range("J12").Value = 816.1667
$ws.Range("L12").Value = 816.1667
$ws.Range("N12").Value = -1156.1667
$ws.Range("H92").Value = 1865.8
$ws.Range("I92").Value = 1930.4546
$ws.Range("J92").Value = 1688
$ws.Range("K92").Value = 1930.4546
$ws.Range("L92").Value = 1688
$ws.Range("M92").Value = -682.4546
$ws.Range("N92").Value = -4184
$ws.Range("H112").Value = 2294.9524
$ws.Range("I112").Value = 1593.25
$ws.Range("K112").Value = 4779.75
$ws.Range("M112").Value = -3671.75
$ws.Range("H116").Value = 10000
$ws.Range("J116").Value = 10000
$ws.Range("L116").Value = 10000
$ws.Range("N116").Value = -16884
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5020
$ws.Range("J32").Value = 4485.3335
$ws.Range("L32").Value = 4485.3335
$ws.Range("N32").Value = -5059.3335
$ws.Range("H61").Value = 5937.1875
$ws.Range("I61").Value = 5999.643
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 5999.643
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -5787.643
$ws.Range("N61").Value = -5924
$ws.Range("H136").Value = 5937.1875
$ws.Range("I136").Value = 5999.643
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 17998.929
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -15448.929
$ws.Range("N136").Value = -21600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1557.4615
$ws.Range("I20").Value = 841.6667
$ws.Range("K20").Value = 841.6667
$ws.Range("M20").Value = -594.6667
$ws.Range("H94").Value = 3250
$ws.Range("I94").Value = 3000
$ws.Range("J94").Value = 3375
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 3375
$ws.Range("M94").Value = -2549
$ws.Range("N94").Value = -4277

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 500
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 500
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -778
$ws.Range("H14").Value = 130
$ws.Range("J14").Value = 175
$ws.Range("L14").Value = 175
$ws.Range("N14").Value = -515
$ws.Range("H21").Value = 10000
$ws.Range("J21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10470
$ws.Range("H99").Value = 3259.8
$ws.Range("I99").Value = 3259.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3259.8
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1761.8
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 2498.8
$ws.Range("I122").Value = 2498.8
$ws.Range("K122").Value = 7496.400000000001
$ws.Range("M122").Value = -5046.400000000001
$ws.Range("H125").Value = 49997.668
$ws.Range("J125").Value = 49997.668
$ws.Range("L125").Value = 49997.668
$ws.Range("N125").Value = -54917.668
$ws.Range("H126").Value = 3259.8
$ws.Range("I126").Value = 3259.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9779.400000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7309.400000000001
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1777.5555
$ws.Range("I5").Value = 1063
$ws.Range("K5").Value = 3189
$ws.Range("M5").Value = -3077
$ws.Range("H12").Value = 70
$ws.Range("I12").Value = 70
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 210
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -37
$ws.Range("N12").ClearContents()
$ws.Range("H18").Value = 2059.261
$ws.Range("I18").Value = 750.82355
$ws.Range("K18").Value = 2252.47065
$ws.Range("M18").Value = -2083.47065
$ws.Range("H33").Value = 142.27272
$ws.Range("I33").Value = 52.5
$ws.Range("J33").Value = 250
$ws.Range("K33").Value = 315
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = -32
$ws.Range("N33").Value = -2066
$ws.Range("H34").Value = 3496.9744
$ws.Range("J34").Value = 3751.2
$ws.Range("L34").Value = 11253.6
$ws.Range("N34").Value = -11421.6
$ws.Range("H68").Value = 1265.3334
$ws.Range("I68").Value = 1026.6666
$ws.Range("J68").Value = 1325
$ws.Range("K68").Value = 3079.9998
$ws.Range("L68").Value = 3975
$ws.Range("M68").Value = -2268.9998
$ws.Range("N68").Value = -5597
$ws.Range("H71").Value = 1265.3334
$ws.Range("I71").Value = 1026.6666
$ws.Range("J71").Value = 1325
$ws.Range("K71").Value = 9239.999400000001
$ws.Range("L71").Value = 11925
$ws.Range("M71").Value = -5183.999400000001
$ws.Range("N71").Value = -20037
$ws.Range("H74").Value = 24230.934
$ws.Range("I74").Value = 24128.5
$ws.Range("K74").Value = 72385.5
$ws.Range("M74").Value = -71324.5
$ws.Range("H77").Value = 24230.934
$ws.Range("I77").Value = 24128.5
$ws.Range("K77").Value = 217156.5
$ws.Range("M77").Value = -211852.5
$ws.Range("H87").Value = 5011.6665
$ws.Range("I87").Value = 5011.6665
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 15034.9995
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -13786.9995
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 5011.6665
$ws.Range("I90").Value = 5011.6665
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 45104.9985
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -38864.9985
$ws.Range("N90").ClearContents()
$ws.Range("H113").Value = 2658.818
$ws.Range("I113").Value = 436.75
$ws.Range("J113").Value = 3928.5715
$ws.Range("K113").Value = 1310.25
$ws.Range("L113").Value = 11785.7145
$ws.Range("M113").Value = 859.75
$ws.Range("N113").Value = -16125.7145
$ws.Range("H117").Value = 2126.5293
$ws.Range("I117").Value = 3594.25
$ws.Range("J117").Value = 1674.9231
$ws.Range("K117").Value = 10782.75
$ws.Range("L117").Value = 5024.7693
$ws.Range("M117").Value = -7340.75
$ws.Range("N117").Value = -11908.7693
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("H135").Value = 1777.5555
$ws.Range("I135").Value = 1063
$ws.Range("K135").Value = 9567
$ws.Range("M135").Value = -7032
$ws.Range("H139").Value = 5663.0557
$ws.Range("I139").Value = 8497
$ws.Range("K139").Value = 25491
$ws.Range("M139").Value = -20351

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6048.067
$ws.Range("I80").Value = 5103.7144
$ws.Range("K80").Value = 5103.7144
$ws.Range("M80").Value = -4105.7144
$ws.Range("H83").Value = 6048.067
$ws.Range("I83").Value = 5103.7144
$ws.Range("K83").Value = 25518.572
$ws.Range("M83").Value = -20526.572
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 35480
$ws.Range("J136").Value = 35480
$ws.Range("L136").Value = 106440
$ws.Range("N136").Value = -111540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2356.9285
$ws.Range("I22").Value = 2356.9285
$ws.Range("K22").Value = 2356.9285
$ws.Range("M22").Value = -2061.9285
$ws.Range("H27").Value = 2356.9285
$ws.Range("I27").Value = 2356.9285
$ws.Range("K27").Value = 2356.9285
$ws.Range("M27").Value = -2249.9285
$ws.Range("H68").Value = 4857.0884
$ws.Range("I68").Value = 3282.35
$ws.Range("K68").Value = 3282.35
$ws.Range("M68").Value = -2533.35
$ws.Range("H71").Value = 4857.0884
$ws.Range("I71").Value = 3282.35
$ws.Range("K71").Value = 16411.75
$ws.Range("M71").Value = -12667.75
$ws.Range("H82").Value = 1897.5385
$ws.Range("I82").Value = 796.5789
$ws.Range("J82").Value = 4885.857
$ws.Range("K82").Value = 796.5789
$ws.Range("L82").Value = 4885.857
$ws.Range("M82").Value = -435.5789
$ws.Range("N82").Value = -5607.857
$ws.Range("H85").Value = 1897.5385
$ws.Range("I85").Value = 796.5789
$ws.Range("J85").Value = 4885.857
$ws.Range("K85").Value = 796.5789
$ws.Range("L85").Value = 4885.857
$ws.Range("M85").Value = 451.4211
$ws.Range("N85").Value = -7381.857
$ws.Range("H93").Value = 3352.5588
$ws.Range("I93").Value = 951.4
$ws.Range("K93").Value = 951.4
$ws.Range("M93").Value = 296.6
$ws.Range("H132").Value = 5384.778
$ws.Range("I132").Value = 3373.5
$ws.Range("J132").Value = 6993.8
$ws.Range("K132").Value = 10120.5
$ws.Range("L132").Value = 20981.4
$ws.Range("M132").Value = -7590.5
$ws.Range("N132").Value = -26041.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11357
$ws.Range("J41").Value = 9612.5
$ws.Range("L41").Value = 9612.5
$ws.Range("N41").Value = -10392.5
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H113").Value = 1442.5
$ws.Range("I113").Value = 1162.8182
$ws.Range("K113").Value = 3488.4546
$ws.Range("M113").Value = -1318.4546
$ws.Range("H136").Value = 2095.8096
$ws.Range("I136").Value = 1312.625
$ws.Range("K136").Value = 3937.875
$ws.Range("M136").Value = -1387.875
